$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, shifting existing rows 39-73 down to 40-74.
$ws.Rows("39").Insert()

# Populate the newly inserted row 39 with the new weekly data entry.
$ws.Cells.Item(39, 1).Value = 10
$ws.Cells.Item(39, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(39, 3).Value = "La Araucanía"
$ws.Cells.Item(39, 4).Value = 44902
$ws.Cells.Item(39, 5).Value = 9
$ws.Cells.Item(39, 6).Value = "Fruta"
$ws.Cells.Item(39, 7).Value = 100103
$ws.Cells.Item(39, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(39, 9).Value = 100103003
$ws.Cells.Item(39, 10).Value = "Damasco"
$ws.Cells.Item(39, 11).Value = "Castle Brite"
$ws.Cells.Item(39, 12).Value = "Primera"
$ws.Cells.Item(39, 13).Value = 200
$ws.Cells.Item(39, 14).Value = 24000
$ws.Cells.Item(39, 15).Value = 24000
$ws.Cells.Item(39, 16).Value = 24000
$ws.Cells.Item(39, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(39, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(39, 19).Value = 1333
$ws.Cells.Item(39, 20).Value = 18
